# Helsinki VRP-ADV version 2 / test_input_3.xlsx update
# Commit: "Add files via upload"
#
# Semantic changes applied:
#  1. "robot" sheet: header in C1 renamed from
#     "average service time per parcel (seconds)" to
#     "service rate (per minute)".
#  2. "problem_input" sheet: new column E added with header
#     "customer arrival rate (per minute)" (left-aligned, same font/
#     border as the other header cells) and values for rows 2-7.

$wb = $excel.ActiveWorkbook

# --- 1. robot sheet: rename the header text in C1 -----------------------
$wsRobot = $wb.Worksheets.Item("robot")
$wsRobot.Range("C1").Value = "service rate (per minute)"

# --- 2. problem_input sheet: add "customer arrival rate" column ---------
$wsInput = $wb.Worksheets.Item("problem_input")

# Copy the header formatting from an existing header cell (bold font +
# thin box border), then switch the new header to left alignment like
# the rest of the edit.
$wsInput.Range("B1").Copy() | Out-Null
$wsInput.Range("E1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wsInput.Range("E1").HorizontalAlignment = -4131      # xlLeft

$wsInput.Range("E1").Value = "customer arrival rate (per minute)"
$wsInput.Range("E2").Value = 0
$wsInput.Range("E3").Value = 10
$wsInput.Range("E4").Value = 5
$wsInput.Range("E5").Value = 10
$wsInput.Range("E6").Value = 5
$wsInput.Range("E7").Value = 5

# Mirror the author's final active selection (cell E6 on problem_input).
$wsInput.Range("E6").Select() | Out-Null
